# Scheduled-runner refresh of market-price / profit columns (H:N) across
# several Leve sheets (currentAveragePrice*, LevePrice*, LeveProfit*).
# Values below were sourced from the upstream data refresh; only the
# already-populated numeric cells in columns H-N are touched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(54, 8).Value = 8750.5
$ws.Cells.Item(54, 9).Value = 2500
$ws.Cells.Item(54, 10).Value = 15001
$ws.Cells.Item(54, 11).Value = 2500
$ws.Cells.Item(54, 12).Value = 15001
$ws.Cells.Item(54, 13).Value = -2014
$ws.Cells.Item(54, 14).Value = -15973

$ws.Cells.Item(81, 8).Value = 30000
$ws.Cells.Item(81, 10).Value = 30000
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 14).Value = -31996

$ws.Cells.Item(84, 8).Value = 30000
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 14).Value = -99984

$ws.Cells.Item(112, 8).Value = 1578.2174
$ws.Cells.Item(112, 10).Value = 2024.9166
$ws.Cells.Item(112, 12).Value = 6074.7498
$ws.Cells.Item(112, 14).Value = -8290.7498

$ws.Cells.Item(113, 8).Value = 4651.25
$ws.Cells.Item(113, 9).Value = 5475.8335
$ws.Cells.Item(113, 10).Value = 3826.6667
$ws.Cells.Item(113, 11).Value = 5475.8335
$ws.Cells.Item(113, 12).Value = 3826.6667
$ws.Cells.Item(113, 13).Value = -2221.8335
$ws.Cells.Item(113, 14).Value = -10334.6667

$ws.Cells.Item(116, 8).Value = 113490
$ws.Cells.Item(116, 9).Value = 145172.86
$ws.Cells.Item(116, 10).Value = 2600
$ws.Cells.Item(116, 11).Value = 145172.86
$ws.Cells.Item(116, 12).Value = 2600
$ws.Cells.Item(116, 13).Value = -141730.86
$ws.Cells.Item(116, 14).Value = -9484

$ws.Cells.Item(125, 8).Value = 7382.6665
$ws.Cells.Item(125, 9).Value = 24692.4
$ws.Cells.Item(125, 10).Value = 725.0769
$ws.Cells.Item(125, 11).Value = 222231.6
$ws.Cells.Item(125, 12).Value = 6525.6921
$ws.Cells.Item(125, 13).Value = -219771.6
$ws.Cells.Item(125, 14).Value = -11445.6921

$ws.Cells.Item(132, 8).Value = 2437.923
$ws.Cells.Item(132, 9).Value = 1931.125
$ws.Cells.Item(132, 10).Value = 3868.8823
$ws.Cells.Item(132, 11).Value = 5793.375
$ws.Cells.Item(132, 12).Value = 11606.6469
$ws.Cells.Item(132, 13).Value = -3263.375
$ws.Cells.Item(132, 14).Value = -16666.6469

$ws.Cells.Item(138, 8).Value = 2390.7
$ws.Cells.Item(138, 9).Value = 1346.362
$ws.Cells.Item(138, 10).Value = 3832.8809
$ws.Cells.Item(138, 11).Value = 4039.086
$ws.Cells.Item(138, 12).Value = 11498.6427
$ws.Cells.Item(138, 13).Value = 1100.914
$ws.Cells.Item(138, 14).Value = -21778.6427

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14, 8).Value = 222222720
$ws.Cells.Item(14, 9).Value = 444444450
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 444444450
$ws.Cells.Item(14, 12).Value = 1000
$ws.Cells.Item(14, 13).Value = -444444275
$ws.Cells.Item(14, 14).Value = -1350

$ws.Cells.Item(61, 8).Value = 265141
$ws.Cells.Item(61, 9).Value = 1887.15
$ws.Cells.Item(61, 10).Value = 557645.25
$ws.Cells.Item(61, 11).Value = 1887.15
$ws.Cells.Item(61, 12).Value = 557645.25
$ws.Cells.Item(61, 13).Value = -1675.15
$ws.Cells.Item(61, 14).Value = -558069.25

$ws.Cells.Item(74, 8).Value = 3247.111
$ws.Cells.Item(74, 9).Value = 907.9459000000001
$ws.Cells.Item(74, 10).Value = 14065.75
$ws.Cells.Item(74, 11).Value = 907.9459000000001
$ws.Cells.Item(74, 12).Value = 14065.75
$ws.Cells.Item(74, 13).Value = -33.94590000000005
$ws.Cells.Item(74, 14).Value = -15813.75

$ws.Cells.Item(77, 8).Value = 3247.111
$ws.Cells.Item(77, 9).Value = 907.9459000000001
$ws.Cells.Item(77, 10).Value = 14065.75
$ws.Cells.Item(77, 11).Value = 4539.7295
$ws.Cells.Item(77, 12).Value = 70328.75
$ws.Cells.Item(77, 13).Value = -171.7295000000004
$ws.Cells.Item(77, 14).Value = -79064.75

$ws.Cells.Item(136, 8).Value = 265141
$ws.Cells.Item(136, 9).Value = 1887.15
$ws.Cells.Item(136, 10).Value = 557645.25
$ws.Cells.Item(136, 11).Value = 5661.450000000001
$ws.Cells.Item(136, 12).Value = 1672935.75
$ws.Cells.Item(136, 13).Value = -3111.450000000001
$ws.Cells.Item(136, 14).Value = -1678035.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(106, 8).Value = 33000.4
$ws.Cells.Item(106, 10).Value = 33000.4
$ws.Cells.Item(106, 12).Value = 33000.4
$ws.Cells.Item(106, 14).Value = -35524.4

$ws.Cells.Item(134, 8).Value = 1727.4814
$ws.Cells.Item(134, 9).Value = 976
$ws.Cells.Item(134, 10).Value = 3230.4443
$ws.Cells.Item(134, 11).Value = 2928
$ws.Cells.Item(134, 12).Value = 9691.332900000001
$ws.Cells.Item(134, 13).Value = -393
$ws.Cells.Item(134, 14).Value = -14761.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 2719
$ws.Cells.Item(15, 10).Value = 2719
$ws.Cells.Item(15, 12).Value = 2719
$ws.Cells.Item(15, 14).Value = -3059

$ws.Cells.Item(31, 8).Value = 16415429
$ws.Cells.Item(31, 9).Value = 166667780
$ws.Cells.Item(31, 10).Value = 24262.727
$ws.Cells.Item(31, 11).Value = 166667780
$ws.Cells.Item(31, 12).Value = 24262.727
$ws.Cells.Item(31, 13).Value = -166667485
$ws.Cells.Item(31, 14).Value = -24852.727

$ws.Cells.Item(34, 8).Value = 16415429
$ws.Cells.Item(34, 9).Value = 166667780
$ws.Cells.Item(34, 10).Value = 24262.727
$ws.Cells.Item(34, 11).Value = 166667780
$ws.Cells.Item(34, 12).Value = 24262.727
$ws.Cells.Item(34, 13).Value = -166667578
$ws.Cells.Item(34, 14).Value = -24666.727

$ws.Cells.Item(93, 8).Value = 6354.857
$ws.Cells.Item(93, 9).Value = 6354.857
$ws.Cells.Item(93, 11).Value = 6354.857
$ws.Cells.Item(93, 13).Value = -4482.857

$ws.Cells.Item(134, 8).Value = 3091.3438
$ws.Cells.Item(134, 9).Value = 3924.4443
$ws.Cells.Item(134, 10).Value = 2020.2142
$ws.Cells.Item(134, 11).Value = 11773.3329
$ws.Cells.Item(134, 12).Value = 6060.642599999999
$ws.Cells.Item(134, 13).Value = -9238.332900000001
$ws.Cells.Item(134, 14).Value = -11130.6426

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 993.3333
$ws.Cells.Item(117, 10).Value = 993.3333
$ws.Cells.Item(117, 12).Value = 2979.9999
$ws.Cells.Item(117, 14).Value = -9863.999899999999

$ws.Cells.Item(122, 8).Value = 1111
$ws.Cells.Item(122, 9).Value = 515.7857
$ws.Cells.Item(122, 10).Value = 1631.8125
$ws.Cells.Item(122, 11).Value = 4642.071300000001
$ws.Cells.Item(122, 12).Value = 14686.3125
$ws.Cells.Item(122, 13).Value = -2192.071300000001
$ws.Cells.Item(122, 14).Value = -19586.3125

$ws.Cells.Item(129, 8).Value = 1318.5714
$ws.Cells.Item(129, 9).Value = 1576.6666
$ws.Cells.Item(129, 10).Value = 1125
$ws.Cells.Item(129, 11).Value = 4729.9998
$ws.Cells.Item(129, 12).Value = 3375
$ws.Cells.Item(129, 13).Value = 270.0002000000004
$ws.Cells.Item(129, 14).Value = -13375

$ws.Cells.Item(131, 8).Value = 891.6667
$ws.Cells.Item(131, 9).Value = 260
$ws.Cells.Item(131, 10).Value = 1035.2273
$ws.Cells.Item(131, 11).Value = 780
$ws.Cells.Item(131, 12).Value = 3105.6819
$ws.Cells.Item(131, 13).Value = 4260
$ws.Cells.Item(131, 14).Value = -13185.6819

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 127700.5
$ws.Cells.Item(122, 9).Value = 169150.67
$ws.Cells.Item(122, 10).Value = 3350
$ws.Cells.Item(122, 11).Value = 507452.01
$ws.Cells.Item(122, 12).Value = 10050
$ws.Cells.Item(122, 13).Value = -505002.01
$ws.Cells.Item(122, 14).Value = -14950

$ws.Cells.Item(132, 8).Value = 4688.5366
$ws.Cells.Item(132, 9).Value = 5462.9287
$ws.Cells.Item(132, 10).Value = 3020.6155
$ws.Cells.Item(132, 11).Value = 16388.7861
$ws.Cells.Item(132, 12).Value = 9061.8465
$ws.Cells.Item(132, 13).Value = -13858.7861
$ws.Cells.Item(132, 14).Value = -14121.8465

$ws.Cells.Item(136, 8).Value = 8453
$ws.Cells.Item(136, 9).Value = 3206.818
$ws.Cells.Item(136, 10).Value = 15666.5
$ws.Cells.Item(136, 11).Value = 9620.454000000002
$ws.Cells.Item(136, 12).Value = 46999.5
$ws.Cells.Item(136, 13).Value = -7070.454000000002
$ws.Cells.Item(136, 14).Value = -52099.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 4000
$ws.Cells.Item(15, 10).Value = 4000
$ws.Cells.Item(15, 12).Value = 4000
$ws.Cells.Item(15, 14).Value = -4576

$ws.Cells.Item(54, 8).Value = 19090.908
$ws.Cells.Item(54, 10).Value = 19090.908
$ws.Cells.Item(54, 12).Value = 19090.908
$ws.Cells.Item(54, 14).Value = -20130.908

$ws.Cells.Item(80, 8).Value = 42999.8
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 42999.8
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 42999.8
$ws.Cells.Item(80, 14).Value = -44995.8
$ws.Cells.Item(80, 13).ClearContents()

$ws.Cells.Item(83, 8).Value = 42999.8
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 42999.8
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 128999.4
$ws.Cells.Item(83, 14).Value = -138983.4
$ws.Cells.Item(83, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 786.7
$ws.Cells.Item(122, 9).Value = 721.1053000000001
$ws.Cells.Item(122, 10).Value = 900
$ws.Cells.Item(122, 11).Value = 2163.3159
$ws.Cells.Item(122, 12).Value = 2700
$ws.Cells.Item(122, 13).Value = 286.6840999999999
$ws.Cells.Item(122, 14).Value = -7600

$ws.Cells.Item(132, 8).Value = 3612.6833
$ws.Cells.Item(132, 9).Value = 4573.154
$ws.Cells.Item(132, 10).Value = 1828.9524
$ws.Cells.Item(132, 11).Value = 13719.462
$ws.Cells.Item(132, 12).Value = 5486.857199999999
$ws.Cells.Item(132, 13).Value = -11189.462
$ws.Cells.Item(132, 14).Value = -10546.8572

Write-Output "Updated 34 leve rows across 7 sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR)."
